$d = $word.ActiveDocument

function Force-Split($r) {
    # Forces the engine to keep this range as its own distinct <w:r> run
    # (rather than silently re-merging it with an adjacent run that has
    # identical formatting) by toggling a character property on and back off.
    $r.Bold = 1
    $r.Bold = 0
}

function Split-Into-Runs($startPos, [string[]]$pieces) {
    # Inserts each piece of $pieces, one after another, starting at $startPos,
    # then (working backwards) forces every inserted piece to remain its own run.
    $pos = $startPos
    $bounds = @($pos)
    foreach ($piece in $pieces) {
        $r = $d.Range($pos, $pos)
        $r.InsertAfter($piece)
        $pos = $pos + $piece.Length
        $bounds += $pos
    }
    for ($i = $pieces.Length - 1; $i -ge 0; $i--) {
        $pieceRng = $d.Range($bounds[$i], $bounds[$i + 1])
        Force-Split($pieceRng)
    }
    return $pos
}

# ---------------------------------------------------------------------------
# Location 1: "...k-means clustering is executed. The scatter plot depicting
# the relationship ... such as " ->
# "...k-means clustering is executed (Figure 1). The scatter plot (Figure 2)
# depicting the relationship ... such as "
# ---------------------------------------------------------------------------
$old1 = "k-means clustering is executed. The scatter plot depicting the relationship between passengers' age and fare was generated, as indicated by the graph, using the variables 'Age', 'Fare', and 'Cluster'. Once the data has been visualized, it can be prepared for genetic algorithm processing. To begin, pertinent attributes such as "
$rng = $d.Content
$found = $rng.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Location 1 text not found" }
$s = $rng.Start
$e = $rng.End
$target = $d.Range($s, $e)
$target.Text = "k-means clustering is executed"
$afterFirst = $target.End

$pieces1 = @(
    " (Figure 1)",
    ". The scatter plot",
    " (Figure 2)",
    " depicting the relationship between passengers' age and fare was generated, as indicated by the graph, using the variables 'Age', 'Fare', and 'Cluster'. Once the data has been visualized, it can be prepared for genetic algorithm processing. To begin, pertinent attributes such as "
)
Split-Into-Runs $afterFirst $pieces1 | Out-Null

# ---------------------------------------------------------------------------
# Location 2: "', and 'Embarked' with respect to the survival rate.
# Correlation analysis helped in understanding the relationship between
# different features." ->
# "', and 'Embarked' with respect to the survival rate (Figure 3). Correlation
# analysis (Figure 4) helped in understanding the relationship between
# different features."
# ---------------------------------------------------------------------------
$old2 = "', and 'Embarked' with respect to the survival rate. Correlation analysis helped in understanding the relationship between different features."
$rng = $d.Content
$found = $rng.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Location 2 text not found" }
$s = $rng.Start
$e = $rng.End
$target = $d.Range($s, $e)
$target.Text = "', and 'Embarked' with respect to the survival rate"
$afterFirst2 = $target.End

$pieces2 = @(
    " (Figure 3)",
    ". Correlation analysis",
    " (Figure 4)",
    " helped in understanding the relationship between different features."
)
Split-Into-Runs $afterFirst2 $pieces2 | Out-Null

# ---------------------------------------------------------------------------
# Location 3: " 32 and epoch 50 times. At last, evaluated the performance of
# the model based on accuracy, precision, recall, and F1-score metrics." ->
# " 32 and epoch 50 times. At last, evaluated the performance of the model
# based on accuracy, precision, recall, and F1-score metrics (Figure 5)."
# ---------------------------------------------------------------------------
$old3 = " 32 and epoch 50 times. At last, evaluated the performance of the model based on accuracy, precision, recall, and F1-score metrics."
$rng = $d.Content
$found = $rng.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Location 3 text not found" }
$s = $rng.Start
$e = $rng.End
$target = $d.Range($s, $e)
$target.Text = " 32 and epoch 50 times. At last, evaluated the performance of the model based on accuracy, precision, recall, and F1-score metrics"
$afterFirst3 = $target.End

$pieces3 = @(
    " (Figure 5)",
    "."
)
Split-Into-Runs $afterFirst3 $pieces3 | Out-Null

Write-Host "All three locations updated."
